$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 42: becomes the "Homo sapiens" entry (was "unassigned")
$ws.Range("A42").Value = "307c55294ffe3b8aa46fce358d55590e"
$ws.Range("B42").Value = "Homo sapiens"
$ws.Range("C42").Value = "Human"
$ws.Range("D42").Value = "Human"
$ws.Range("J42").ClearContents()

# Row 43: becomes the "unassigned" entry (was "Homo sapiens")
$ws.Range("A43").Value = "c0a3f3ed23f04247d92740a9502f8b57"
$ws.Range("B43").Value = "unassigned"
$ws.Range("C43").Value = "unassigned"
$ws.Range("D43").Value = "unassigned"
$ws.Range("J43").Value = 0

# Row 55: becomes the "unassigned" entry (was "Centropristis striata")
$ws.Range("A55").Value = "5e733a21f67e541f28ed4bf4fe025044"
$ws.Range("B55").Value = "unassigned"
$ws.Range("C55").Value = "unassigned"
$ws.Range("D55").Value = "unassigned"

# Row 56: becomes the "Centropristis striata" entry (was "unassigned")
$ws.Range("A56").Value = "975b1dbdc7405f6e27bf63893e91e0ed"
$ws.Range("B56").Value = "Centropristis striata"
$ws.Range("C56").Value = "Black sea bass"
$ws.Range("D56").Value = "Teleost Fish"
